$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1956.5454
$ws.Range("J17").Value = 1956.5454
$ws.Range("L17").Value = 5869.6362
$ws.Range("N17").Value = -6205.6362
$ws.Range("H87").Value = 67231.375
$ws.Range("J87").Value = 86770.2
$ws.Range("L87").Value = 86770.2
$ws.Range("N87").Value = -89266.2
$ws.Range("H90").Value = 67231.375
$ws.Range("J90").Value = 86770.2
$ws.Range("L90").Value = 260310.6
$ws.Range("N90").Value = -272790.6
$ws.Range("H93").Value = 50000.0
$ws.Range("J93").Value = 50000.0
$ws.Range("L93").Value = 50000.0
$ws.Range("N93").Value = -54992.0
$ws.Range("H106").Value = 2091.8462
$ws.Range("I106").Value = 1529.4
$ws.Range("K106").Value = 1529.4
$ws.Range("M106").Value = -898.4000000000001
$ws.Range("H112").Value = 1455.7037
$ws.Range("J112").Value = 1503.76
$ws.Range("L112").Value = 4511.28
$ws.Range("N112").Value = -6727.28
$ws.Range("H127").Value = 4595.6
$ws.Range("I127").Value = 4595.6
$ws.Range("J127").Value = 0.0
$ws.Range("K127").Value = 13786.8
$ws.Range("L127").Value = 0.0
$ws.Range("M127").Value = -8826.800000000001
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 1962.5555
$ws.Range("I132").Value = 1932.0294
$ws.Range("K132").Value = 5796.0882
$ws.Range("M132").Value = -3266.0882
$ws.Range("H133").Value = 79851.664
$ws.Range("J133").Value = 79851.664
$ws.Range("L133").Value = 79851.664
$ws.Range("N133").Value = -89971.664
$ws.Range("H136").Value = 100000.0
$ws.Range("J136").Value = 100000.0
$ws.Range("L136").Value = 100000.0
$ws.Range("N136").Value = -110200.0
$ws.Range("H138").Value = 2164.9092
$ws.Range("I138").Value = 1547.0
$ws.Range("J138").Value = 3359.5334
$ws.Range("K138").Value = 4641.0
$ws.Range("L138").Value = 10078.6002
$ws.Range("M138").Value = 499.0
$ws.Range("N138").Value = -20358.6002
$ws.Range("H141").Value = 1999.5
$ws.Range("J141").Value = 0.0
$ws.Range("L141").Value = 0.0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7847.4585
$ws.Range("I32").Value = 9463.611
$ws.Range("K32").Value = 9463.611
$ws.Range("M32").Value = -9176.611
$ws.Range("H102").Value = 3046.2632
$ws.Range("I102").Value = 2598.6
$ws.Range("K102").Value = 2598.6
$ws.Range("M102").Value = -976.5999999999999
$ws.Range("H110").Value = 574.72
$ws.Range("I110").Value = 587.087
$ws.Range("K110").Value = 587.087
$ws.Range("M110").Value = 1457.913
$ws.Range("H122").Value = 2883.149
$ws.Range("I122").Value = 2478.2646
$ws.Range("J122").Value = 3942.077
$ws.Range("K122").Value = 7434.793799999999
$ws.Range("L122").Value = 11826.231
$ws.Range("M122").Value = -4984.793799999999
$ws.Range("N122").Value = -16726.231
$ws.Range("H134").Value = 0.0
$ws.Range("J134").Value = 0.0
$ws.Range("L134").Value = 0.0
$ws.Range("N134").ClearContents()
$ws.Range("H139").Value = 88570.71
$ws.Range("J139").Value = 88570.71
$ws.Range("L139").Value = 88570.71
$ws.Range("N139").Value = -98850.71

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 0.0
$ws.Range("J110").Value = 0.0
$ws.Range("L110").Value = 0.0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 66583.84
$ws.Range("J68").Value = 68049.164
$ws.Range("L68").Value = 68049.164
$ws.Range("N68").Value = -69547.164
$ws.Range("H71").Value = 66583.84
$ws.Range("J71").Value = 68049.164
$ws.Range("L71").Value = 204147.492
$ws.Range("N71").Value = -211635.492
$ws.Range("H86").Value = 5700.0
$ws.Range("I86").Value = 5000.0
$ws.Range("K86").Value = 5000.0
$ws.Range("M86").Value = -3877.0
$ws.Range("H89").Value = 5700.0
$ws.Range("I89").Value = 5000.0
$ws.Range("K89").Value = 25000.0
$ws.Range("M89").Value = -19384.0
$ws.Range("H122").Value = 4099.625
$ws.Range("I122").Value = 1999.25
$ws.Range("J122").Value = 6200.0
$ws.Range("K122").Value = 5997.75
$ws.Range("L122").Value = 18600.0
$ws.Range("M122").Value = -3547.75
$ws.Range("N122").Value = -23500.0
$ws.Range("H132").Value = 2064.25
$ws.Range("I132").Value = 2419.0
$ws.Range("J132").Value = 1000.0
$ws.Range("K132").Value = 7257.0
$ws.Range("L132").Value = 3000.0
$ws.Range("M132").Value = -4727.0
$ws.Range("N132").Value = -8060.0
$ws.Range("H134").Value = 1588.4839
$ws.Range("I134").Value = 1588.4839
$ws.Range("K134").Value = 4765.4517
$ws.Range("M134").Value = -2230.4517

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 90000.0
$ws.Range("J37").Value = 90000.0
$ws.Range("L37").Value = 270000.0
$ws.Range("N37").Value = -270224.0
$ws.Range("H60").Value = 1743.6666
$ws.Range("I60").Value = 1447.6364
$ws.Range("K60").Value = 4342.9092
$ws.Range("M60").Value = -4091.9092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 28928.5
$ws.Range("J47").Value = 28928.5
$ws.Range("L47").Value = 28928.5
$ws.Range("N47").Value = -30064.5
$ws.Range("H55").Value = 7336.1665
$ws.Range("J55").Value = 7750.8
$ws.Range("L55").Value = 7750.8
$ws.Range("N55").Value = -8404.8
$ws.Range("H111").Value = 0.0
$ws.Range("J111").Value = 0.0
$ws.Range("L111").Value = 0.0
$ws.Range("N111").ClearContents()
$ws.Range("H126").Value = 3793.8147
$ws.Range("J126").Value = 4094.762
$ws.Range("L126").Value = 12284.286
$ws.Range("N126").Value = -17224.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4733.25
$ws.Range("I68").Value = 3000.0
$ws.Range("J68").Value = 5311.0
$ws.Range("K68").Value = 3000.0
$ws.Range("L68").Value = 5311.0
$ws.Range("M68").Value = -2251.0
$ws.Range("N68").Value = -6809.0
$ws.Range("H71").Value = 4733.25
$ws.Range("I71").Value = 3000.0
$ws.Range("J71").Value = 5311.0
$ws.Range("K71").Value = 15000.0
$ws.Range("L71").Value = 26555.0
$ws.Range("M71").Value = -11256.0
$ws.Range("N71").Value = -34043.0
$ws.Range("H104").Value = 95691.0
$ws.Range("J104").Value = 95691.0
$ws.Range("L104").Value = 95691.0
$ws.Range("N104").Value = -102679.0
$ws.Range("H110").Value = 80912.664
$ws.Range("J110").Value = 80912.664
$ws.Range("L110").Value = 80912.664
$ws.Range("N110").Value = -89092.664
$ws.Range("H136").Value = 7491.32
$ws.Range("I136").Value = 6752.421
$ws.Range("J136").Value = 9831.167
$ws.Range("K136").Value = 20257.263
$ws.Range("L136").Value = 29493.501
$ws.Range("M136").Value = -17707.263
$ws.Range("N136").Value = -34593.501
$ws.Range("H141").Value = 0.0
$ws.Range("J141").Value = 0.0
$ws.Range("L141").Value = 0.0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 0.0
$ws.Range("J104").Value = 0.0
$ws.Range("L104").Value = 0.0
$ws.Range("N104").ClearContents()
$ws.Range("H107").Value = 2116.1904
$ws.Range("I107").Value = 817.5
$ws.Range("K107").Value = 2452.5
$ws.Range("M107").Value = -532.5
$ws.Range("H114").Value = 80118.4
$ws.Range("J114").Value = 80118.4
$ws.Range("L114").Value = 80118.4
$ws.Range("N114").Value = -88796.4
$ws.Range("H126").Value = 4798.5
$ws.Range("I126").Value = 4161.8184
$ws.Range("J126").Value = 7133.0
$ws.Range("K126").Value = 12485.4552
$ws.Range("L126").Value = 21399.0
$ws.Range("M126").Value = -10015.4552
$ws.Range("N126").Value = -26339.0
$ws.Range("H132").Value = 5462.3687
$ws.Range("I132").Value = 4710.5
$ws.Range("J132").Value = 8281.875
$ws.Range("K132").Value = 14131.5
$ws.Range("L132").Value = 24845.625
$ws.Range("M132").Value = -11601.5
$ws.Range("N132").Value = -29905.625
$ws.Range("H136").Value = 14403030.0
$ws.Range("I136").Value = 20331372.0
$ws.Range("K136").Value = 60994116.0
$ws.Range("M136").Value = -60991566.0

